$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATOS")

# Rearrange columns: A=SKU, B=DESCRIPCION, C=PIEZAS, D=PRECIO
$ws.Range("A1").Value = "SKU"
$ws.Range("B1").Value = "DESCRIPCION"
$ws.Range("C1").Value = "PIEZAS"
$ws.Range("D1").Value = "PRECIO"

$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 1234567890
$ws.Range("C2").Value = "Perrito"
$ws.Range("D2").Value = 3

$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 1213487942
$ws.Range("C3").Value = "Gatito"
$ws.Range("D3").Value = 4

$ws.Activate()
$excel.ActiveWindow.Zoom = 160
$ws.Range("D5").Select() | Out-Null

Write-Output "done"
